$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update simple text/number values
$ws.Range("A3").Value = 7709182030
$ws.Range("A4").Value = "MBA (ITSM)"
$ws.Range("A5").Value = 77122558691
$ws.Range("A6").Value = "Validity End: 30-Jun-2026"

# Turn A2 into a hyperlink pointing at the new email address
$ws.Range("A2").Value = "darshanchawade@gmail.com"
$ws.Hyperlinks.Add($ws.Range("A2"), "mailto:darshanchawade@gmail.com")

# Update the selection to match the target workbook (A10 active cell)
$ws.Range("A10").Select()
